$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Notes" column header
$ws.Range("M1").Value = "Notes"

# Add note about planned student absence for an exam in another course
$ws.Range("M9").Value = "EH exam in other course"
$ws.Range("M15").Value = "EH exam in other course"

# Update the selected cell to match the final state in the file
$ws.Range("B15").Select()
